$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (cols that actually change width) ---
# Target stored widths (per XML) are on a 1/256 "character" grid, but this
# host's ColumnWidth setter snaps to an MDW-7 pixel grid (1/6 steps), so we
# feed it the ColumnWidth input whose rounded result lands nearest the
# intended stored width.
$ws.Columns.Item(6).ColumnWidth  = 1.3333333333333333   # F: 3.140625 -> ~2.140625
$ws.Columns.Item(7).ColumnWidth  = 1.3333333333333333   # G: 3.140625 -> ~2.140625
$ws.Columns.Item(8).ColumnWidth  = 2.3333333333333335   # H: 2.140625 -> ~3.140625
$ws.Columns.Item(9).ColumnWidth  = 2.3333333333333335   # I: 2.140625 -> ~3.140625
$ws.Columns.Item(10).ColumnWidth = 2.3333333333333335   # J: 2.140625 -> ~3.140625
$ws.Columns.Item(13).ColumnWidth = 6.833333333333333    # M: 5.7109375 -> ~7.7109375
$ws.Columns.Item(16).ColumnWidth = 3.8333333333333335   # P: 5.7109375 -> ~4.7109375
$ws.Columns.Item(17).ColumnWidth = 3.8333333333333335   # Q: 5.7109375 -> ~4.7109375

# --- Row 1 cell value changes ---
$ws.Range("C1").Value = 30
$ws.Range("D1").Value = 11
$ws.Range("E1").Value = 20
$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 9
$ws.Range("H1").Value = 33
$ws.Range("I1").Value = 25
$ws.Range("J1").Value = 29
$ws.Range("K1").Value = 19
$ws.Range("L1").Value = 17
$ws.Range("M1").Value = 0.021989999999999999
$ws.Range("N1").Value = 0.060999999999999999
$ws.Range("O1").Value = 0.002
$ws.Range("P1").Value = 0.040000000000000001
$ws.Range("Q1").Value = 0.029999999999999999
